$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Year 1 (P.I. Summer) / Year 3 (Post-Doc) escalation toggle swap ---
# Row 8: move the "switched on" escalation multiplier from column D to column C.
$ws.Range("C8").Formula = "=`$K8*1.03*1.03*10%*1"
$ws.Range("D8").Formula = "=`$K8*1.03*1.03*1.03*10%*0"

# Row 10: Post-Doctoral Fellow now spreads pay over 10 months instead of 9.
$ws.Range("B10").Formula = "=K9*(1/10)"

# Senior Personnel - Downtown Campus: base salary reduced.
$ws.Range("K15").Value = 42000

# Graduate Assistant(s) (Enrolled): number of students doubled.
$ws.Range("K22").Value = 2

# Undergraduate Assistant(s) (Unenrolled - summer): hours/week reduced.
$ws.Range("L23").Value = 10

# --- View state: scroll position / active selection ---
$ws.Range("B1").Select()
$ws.Range("C9").Select()
